$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new display text. Values are written as plain
# text (matching the inline-string storage already used in this sheet)
# by temporarily forcing a text number format, assigning the value, then
# clearing formats again so the cell keeps its original (default) style.
$updates = [ordered]@{
    "D2" = "305.62"
    "E2" = "1.18%"
    "D3" = "36.01"
    "E3" = "-1.60%"
    "D4" = "5.064"
    "E4" = "1.62%"
    "D5" = "0.07999"
    "E5" = "3.14%"
    "D6" = "2.186"
    "E6" = "4.06%"
    "D7" = "8.010"
    "E7" = "1.28%"
    "D8" = "4.151"
    "E8" = "2.87%"
    "D9" = "0.9292"
    "E9" = "0.73%"
    "D10" = "0.09890"
    "E10" = "1.07%"
    "D11" = "0.1871"
    "E11" = "0.05%"
    "D12" = "0.09029"
    "E12" = "4.61%"
    "D13" = "0.03617"
    "E13" = "2.98%"
    "D14" = "0.09911"
    "E14" = "-0.45%"
    "D15" = "0.001460"
    "E15" = "-0.36%"
    "D16" = "0.005656"
    "E16" = "-0.12%"
    "D17" = "3.449"
    "E17" = "-0.45%"
    "E18" = "13.84%"
    "D19" = "0.3372"
    "E19" = "-1.05%"
    "D20" = "0.1355"
    "E20" = "1.00%"
    "D21" = "5.062"
    "E21" = "6.03%"
    "D22" = "0.2190"
    "E22" = "-0.35%"
    "D23" = "0.04592"
    "E23" = "0.02%"
    "D24" = "0.001240"
    "E24" = "0.89%"
    "D25" = "0.004756"
    "E25" = "-6.43%"
    "D26" = "0.0001302"
    "E26" = "-6.90%"
    "D27" = "0.0004501"
    "E27" = "65.13%"
    "D39" = "0.01945"
    "E39" = "10.09%"
    "D40" = "0.04903"
    "E40" = "5.29%"
    "D41" = "0.007764"
    "E41" = "4.25%"
    "E42" = "0.21%"
    "D43" = "0.007804"
    "E43" = "1.31%"
    "D44" = "0.002107"
    "E44" = "-5.83%"
    "D45" = "0.01145"
    "E45" = "10.07%"
    "D46" = "0.00006208"
    "E46" = "0.24%"
    "D47" = "0.00000000750"
    "E47" = "0.24%"
    "D48" = "51.99"
    "E48" = "36.23%"
    "D49" = "0.001801"
    "E49" = "-9.78%"
    "D50" = "0.00002101"
    "E50" = "0.24%"
    "D51" = "0.0002001"
    "E51" = "0.24%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
